$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: the existing '总计' sheet (sheetId 4) becomes '2022-Q1'.
#         It is renamed in place so it keeps sheetId 4, matching the
#         target diff (old totals data is fully overwritten with the
#         2022-Q1 per-fund holdings table).
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# A helper cell with default (unstyled) formatting used as a format
# donor so numeric-looking text (fund codes, percentages, ...) can be
# written as literal text without leaving a stray NumberFormat style
# behind once we reset it.
$q1Plain = $q1.Range("Z100")

# Header row + column-A style already exist on sheet '2021-Q4' with
# exactly the same layout (style index 2: bold, centered, bordered).
# Re-use that as the format donor for the new columns / extra rows.
$fmtSrc = $wb.Worksheets.Item("2021-Q4")
$fmtSrc.Range("E1:H1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2").Copy()
$q1.Range("A2:A23").PasteSpecial(-4122)

# ---- header row ----
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# ---- data rows ----
# row 2
$q1.Range("A2").Value = 0
$c = $q1.Range("B2"); $c.NumberFormat = "@"; $c.Value = "009983"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C2").Value = "永赢港股通品质生活慧选混合"
$c = $q1.Range("D2"); $c.NumberFormat = "@"; $c.Value = "9.75"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E2"); $c.NumberFormat = "@"; $c.Value = "87.97"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F2"); $c.NumberFormat = "@"; $c.Value = "5.04"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G2"); $c.NumberFormat = "@"; $c.Value = "0.4914"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H2").Value = 6
# row 3
$q1.Range("A3").Value = 1
$c = $q1.Range("B3"); $c.NumberFormat = "@"; $c.Value = "003713"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C3").Value = "英大睿盛灵活配置混合A"
$c = $q1.Range("D3"); $c.NumberFormat = "@"; $c.Value = "5.99"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E3"); $c.NumberFormat = "@"; $c.Value = "87.42"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F3"); $c.NumberFormat = "@"; $c.Value = "7.03"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G3"); $c.NumberFormat = "@"; $c.Value = "0.4211"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H3").Value = 2
# row 4
$q1.Range("A4").Value = 2
$c = $q1.Range("B4"); $c.NumberFormat = "@"; $c.Value = "011315"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C4").Value = "永赢港股通优质成长一年持有期混合型证券投资基金"
$c = $q1.Range("D4"); $c.NumberFormat = "@"; $c.Value = "4.07"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E4"); $c.NumberFormat = "@"; $c.Value = "90.36"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F4"); $c.NumberFormat = "@"; $c.Value = "5.01"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G4"); $c.NumberFormat = "@"; $c.Value = "0.2039"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H4").Value = 6
# row 5
$q1.Range("A5").Value = 3
$c = $q1.Range("B5"); $c.NumberFormat = "@"; $c.Value = "008480"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C5").Value = "永赢股息优选混合A"
$c = $q1.Range("D5"); $c.NumberFormat = "@"; $c.Value = "3.42"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E5"); $c.NumberFormat = "@"; $c.Value = "89.21"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F5"); $c.NumberFormat = "@"; $c.Value = "5.45"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G5"); $c.NumberFormat = "@"; $c.Value = "0.1864"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H5").Value = 6
# row 6
$q1.Range("A6").Value = 4
$c = $q1.Range("B6"); $c.NumberFormat = "@"; $c.Value = "003714"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C6").Value = "英大睿盛灵活配置混合C"
$c = $q1.Range("D6"); $c.NumberFormat = "@"; $c.Value = "2.40"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E6"); $c.NumberFormat = "@"; $c.Value = "87.42"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F6"); $c.NumberFormat = "@"; $c.Value = "7.03"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G6"); $c.NumberFormat = "@"; $c.Value = "0.1687"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H6").Value = 2
# row 7
$q1.Range("A7").Value = 5
$c = $q1.Range("B7"); $c.NumberFormat = "@"; $c.Value = "013393"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C7").Value = "信达澳银价值精选混合A"
$c = $q1.Range("D7"); $c.NumberFormat = "@"; $c.Value = "3.61"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E7"); $c.NumberFormat = "@"; $c.Value = "81.31"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F7"); $c.NumberFormat = "@"; $c.Value = "3.49"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G7"); $c.NumberFormat = "@"; $c.Value = "0.1260"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H7").Value = 2
# row 8
$q1.Range("A8").Value = 6
$c = $q1.Range("B8"); $c.NumberFormat = "@"; $c.Value = "005526"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C8").Value = "工银瑞信新生代消费灵活配置混合"
$c = $q1.Range("D8"); $c.NumberFormat = "@"; $c.Value = "2.23"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E8"); $c.NumberFormat = "@"; $c.Value = "94.09"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F8"); $c.NumberFormat = "@"; $c.Value = "4.75"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G8"); $c.NumberFormat = "@"; $c.Value = "0.1059"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H8").Value = 4
# row 9
$q1.Range("A9").Value = 7
$c = $q1.Range("B9"); $c.NumberFormat = "@"; $c.Value = "004995"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C9").Value = "广发品牌消费股票A"
$c = $q1.Range("D9"); $c.NumberFormat = "@"; $c.Value = "2.73"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E9"); $c.NumberFormat = "@"; $c.Value = "92.24"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F9"); $c.NumberFormat = "@"; $c.Value = "3.47"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G9"); $c.NumberFormat = "@"; $c.Value = "0.0947"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H9").Value = 8
# row 10
$q1.Range("A10").Value = 8
$c = $q1.Range("B10"); $c.NumberFormat = "@"; $c.Value = "002005"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C10").Value = "工银瑞信新得利混合"
$c = $q1.Range("D10"); $c.NumberFormat = "@"; $c.Value = "3.82"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E10"); $c.NumberFormat = "@"; $c.Value = "29.93"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F10"); $c.NumberFormat = "@"; $c.Value = "2.18"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G10"); $c.NumberFormat = "@"; $c.Value = "0.0833"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H10").Value = 3
# row 11
$q1.Range("A11").Value = 9
$c = $q1.Range("B11"); $c.NumberFormat = "@"; $c.Value = "009240"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C11").Value = "泰康蓝筹优势一年持有期股票"
$c = $q1.Range("D11"); $c.NumberFormat = "@"; $c.Value = "3.72"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E11"); $c.NumberFormat = "@"; $c.Value = "92.20"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F11"); $c.NumberFormat = "@"; $c.Value = "1.67"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G11"); $c.NumberFormat = "@"; $c.Value = "0.0621"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H11").Value = 10
# row 12
$q1.Range("A12").Value = 10
$c = $q1.Range("B12"); $c.NumberFormat = "@"; $c.Value = "233008"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C12").Value = "大摩消费领航混合基金"
$c = $q1.Range("D12"); $c.NumberFormat = "@"; $c.Value = "0.89"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E12"); $c.NumberFormat = "@"; $c.Value = "79.72"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F12"); $c.NumberFormat = "@"; $c.Value = "6.10"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G12"); $c.NumberFormat = "@"; $c.Value = "0.0543"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H12").Value = 7
# row 13
$q1.Range("A13").Value = 11
$c = $q1.Range("B13"); $c.NumberFormat = "@"; $c.Value = "003446"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C13").Value = "英大睿鑫灵活配置混合A"
$c = $q1.Range("D13"); $c.NumberFormat = "@"; $c.Value = "0.59"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E13"); $c.NumberFormat = "@"; $c.Value = "89.46"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F13"); $c.NumberFormat = "@"; $c.Value = "7.42"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G13"); $c.NumberFormat = "@"; $c.Value = "0.0438"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H13").Value = 1
# row 14
$q1.Range("A14").Value = 12
$c = $q1.Range("B14"); $c.NumberFormat = "@"; $c.Value = "003447"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C14").Value = "英大睿鑫灵活配置混合C"
$c = $q1.Range("D14"); $c.NumberFormat = "@"; $c.Value = "0.51"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E14"); $c.NumberFormat = "@"; $c.Value = "89.46"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F14"); $c.NumberFormat = "@"; $c.Value = "7.42"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G14"); $c.NumberFormat = "@"; $c.Value = "0.0378"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H14").Value = 1
# row 15
$q1.Range("A15").Value = 13
$c = $q1.Range("B15"); $c.NumberFormat = "@"; $c.Value = "004987"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C15").Value = "诺德新享灵活配置混合"
$c = $q1.Range("D15"); $c.NumberFormat = "@"; $c.Value = "0.80"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E15"); $c.NumberFormat = "@"; $c.Value = "88.76"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F15"); $c.NumberFormat = "@"; $c.Value = "4.39"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G15"); $c.NumberFormat = "@"; $c.Value = "0.0351"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H15").Value = 10
# row 16
$q1.Range("A16").Value = 14
$c = $q1.Range("B16"); $c.NumberFormat = "@"; $c.Value = "010245"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C16").Value = "广发品牌消费股票C"
$c = $q1.Range("D16"); $c.NumberFormat = "@"; $c.Value = "0.52"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E16"); $c.NumberFormat = "@"; $c.Value = "92.24"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F16"); $c.NumberFormat = "@"; $c.Value = "3.47"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G16"); $c.NumberFormat = "@"; $c.Value = "0.0180"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H16").Value = 8
# row 17
$q1.Range("A17").Value = 15
$c = $q1.Range("B17"); $c.NumberFormat = "@"; $c.Value = "007133"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C17").Value = "嘉实长青竞争优势股票A"
$c = $q1.Range("D17"); $c.NumberFormat = "@"; $c.Value = "0.25"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E17"); $c.NumberFormat = "@"; $c.Value = "90.07"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F17"); $c.NumberFormat = "@"; $c.Value = "5.43"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G17"); $c.NumberFormat = "@"; $c.Value = "0.0136"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H17").Value = 3
# row 18
$q1.Range("A18").Value = 16
$c = $q1.Range("B18"); $c.NumberFormat = "@"; $c.Value = "000679"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C18").Value = "招商丰利灵活配置混合A"
$c = $q1.Range("D18"); $c.NumberFormat = "@"; $c.Value = "0.39"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E18"); $c.NumberFormat = "@"; $c.Value = "74.75"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F18"); $c.NumberFormat = "@"; $c.Value = "3.44"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G18"); $c.NumberFormat = "@"; $c.Value = "0.0134"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H18").Value = 8
# row 19
$q1.Range("A19").Value = 17
$c = $q1.Range("B19"); $c.NumberFormat = "@"; $c.Value = "013394"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C19").Value = "信达澳银价值精选混合C"
$c = $q1.Range("D19"); $c.NumberFormat = "@"; $c.Value = "0.37"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E19"); $c.NumberFormat = "@"; $c.Value = "81.31"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F19"); $c.NumberFormat = "@"; $c.Value = "3.49"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G19"); $c.NumberFormat = "@"; $c.Value = "0.0129"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H19").Value = 2
# row 20
$q1.Range("A20").Value = 18
$c = $q1.Range("B20"); $c.NumberFormat = "@"; $c.Value = "008481"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C20").Value = "永赢股息优选混合C"
$c = $q1.Range("D20"); $c.NumberFormat = "@"; $c.Value = "0.19"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E20"); $c.NumberFormat = "@"; $c.Value = "89.21"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F20"); $c.NumberFormat = "@"; $c.Value = "5.45"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G20"); $c.NumberFormat = "@"; $c.Value = "0.0104"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H20").Value = 6
# row 21
$q1.Range("A21").Value = 19
$c = $q1.Range("B21"); $c.NumberFormat = "@"; $c.Value = "001608"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C21").Value = "英大策略优选混合C"
$c = $q1.Range("D21"); $c.NumberFormat = "@"; $c.Value = "0.03"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E21"); $c.NumberFormat = "@"; $c.Value = "89.86"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F21"); $c.NumberFormat = "@"; $c.Value = "8.66"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G21"); $c.NumberFormat = "@"; $c.Value = "0.0026"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H21").Value = 1
# row 22
$q1.Range("A22").Value = 20
$c = $q1.Range("B22"); $c.NumberFormat = "@"; $c.Value = "007134"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C22").Value = "嘉实长青竞争优势股票C"
$c = $q1.Range("D22"); $c.NumberFormat = "@"; $c.Value = "0.02"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E22"); $c.NumberFormat = "@"; $c.Value = "90.07"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F22"); $c.NumberFormat = "@"; $c.Value = "5.43"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G22"); $c.NumberFormat = "@"; $c.Value = "0.0011"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H22").Value = 3
# row 23
$q1.Range("A23").Value = 21
$c = $q1.Range("B23"); $c.NumberFormat = "@"; $c.Value = "002416"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("C23").Value = "招商丰利灵活配置混合C"
$c = $q1.Range("D23"); $c.NumberFormat = "@"; $c.Value = "0.02"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("E23"); $c.NumberFormat = "@"; $c.Value = "74.75"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("F23"); $c.NumberFormat = "@"; $c.Value = "3.44"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$c = $q1.Range("G23"); $c.NumberFormat = "@"; $c.Value = "0.0007"; $q1Plain.Copy(); $c.PasteSpecial(-4122)
$q1.Range("H23").Value = 8

# ------------------------------------------------------------------
# Step 2: add a brand-new '总计' sheet (sheetId 5) right after
#         '2022-Q1', matching the target order & sheetId sequence.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Add()
$total.Name = "TotalsStaging"
$wb.Worksheets.Item("TotalsStaging").Move($null, $wb.Worksheets.Item("2022-Q1"))
$wb.Worksheets.Item("TotalsStaging").Name = "总计"
$total = $wb.Worksheets.Item("总计")

# format donor for header/column-A style (index 2), and a plain
# (unstyled) donor cell for text-looking numeric data.
$fmtSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$fmtSrc.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)
$totalPlain = $total.Range("Z100")

# ---- header row ----
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# ---- data rows ----
# row 2
$total.Range("A2").Value = 0
$c = $total.Range("B2"); $c.NumberFormat = "@"; $c.Value = "2022-Q1"; $totalPlain.Copy(); $c.PasteSpecial(-4122)
$total.Range("C2").Value = 22
$total.Range("D2").Value = 2.19
# row 3
$total.Range("A3").Value = 1
$c = $total.Range("B3"); $c.NumberFormat = "@"; $c.Value = "2021-Q4"; $totalPlain.Copy(); $c.PasteSpecial(-4122)
$total.Range("C3").Value = 6
$total.Range("D3").Value = 0.66
# row 4
$total.Range("A4").Value = 2
$c = $total.Range("B4"); $c.NumberFormat = "@"; $c.Value = "2021-Q3"; $totalPlain.Copy(); $c.PasteSpecial(-4122)
$total.Range("C4").Value = 10
$total.Range("D4").Value = 0.54
# row 5
$total.Range("A5").Value = 3
$c = $total.Range("B5"); $c.NumberFormat = "@"; $c.Value = "2020-Q4"; $totalPlain.Copy(); $c.PasteSpecial(-4122)
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.26
